$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that are no longer present in the target (2013-2020 era rows
# 15-22 corresponding to the old sheet layout), shrinking the sheet so the
# used range ends at row 14.
$ws.Rows("15:22").Delete()

# Rewrite the year labels and values for rows 2-14 with the updated dataset.
$years = @("2010年","2011年","2012年","2013年","2014年","2015年","2016年","2017年","2018年","2019年","2020年","2021年","2022年")
$values = @(-5293598.8449153, -1376574.77202481, -8707436.33339946, -6292478.71684688, -6687267.23382289, -20181550.0279276, -21858696.9798611, -20660573.5263438, -17735761.5754534, -12918094.3532083, -15876286.877763, -13446564.4952935, -9055666.28321301)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value() = $years[$i]
    $ws.Cells.Item($row, 2).Value() = $values[$i]
}
